$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "date" column (F2:F7) forward by 13 days (chapter 9/10 merge)
$ws.Range("F2").Value = 44633
$ws.Range("F3").Value = 44632
$ws.Range("F4").Value = 44631
$ws.Range("F5").Value = 44630
$ws.Range("F6").Value = 44629
$ws.Range("F7").Value = 44628
